$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from H1 to I1:J1, then set header text
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill I and J columns for rows 2-72
$iVals = @(9,9,6,9,8,7,7,8,8,7,8,7,7,7,7,6,6,7,6,7,10,8,7,8,5,9,8,9,9,9,9,7,9,9,9,9,9,9,9,8,8,9,9,8,9,9,9,9,10,9,9,9,9,9,9,9,9,7,9,9,9,9,7,9,4,2,6,8,3,5,5)
$jVals = @(9,9,6,9,8,7,7,8,8,7,8,7,7,7,7,6,6,7,6,7,10,8,8,8,5,9,8,9,9,9,9,7,9,9,9,9,9,9,9,8,8,9,9,9,9,9,10,9,10,9,9,9,9,9,9,9,9,7,10,9,9,9,7,9,4,3,6,8,3,5,5)

for ($idx = 0; $idx -lt $iVals.Count; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$idx]
    $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}
